# Update "想去人数" (interest count) values in column F across all four
# sheets to reflect the latest scrape, per commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 146
$ws.Range("F5").Value = 395
$ws.Range("F6").Value = 778
$ws.Range("F7").Value = 224
$ws.Range("F8").Value = 1116
$ws.Range("F9").Value = 298
$ws.Range("F12").Value = 651
$ws.Range("F14").Value = 500
$ws.Range("F18").Value = 2857
$ws.Range("F21").Value = 25
$ws.Range("F24").Value = 219
$ws.Range("F26").Value = 2438
$ws.Range("F28").Value = 974
$ws.Range("F31").Value = 271
$ws.Range("F32").Value = 1056
$ws.Range("F35").Value = 278

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1067
$ws.Range("F5").Value = 1067
$ws.Range("F11").Value = 7
$ws.Range("F15").Value = 595
$ws.Range("F18").Value = 975
$ws.Range("F27").Value = 3839
$ws.Range("F32").Value = 37
$ws.Range("F35").Value = 29

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 75
$ws.Range("F5").Value = 2425
$ws.Range("F6").Value = 1008
$ws.Range("F9").Value = 1269
$ws.Range("F10").Value = 338
$ws.Range("F11").Value = 91

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 2425
$ws.Range("F6").Value = 1008
$ws.Range("F7").Value = 1269
$ws.Range("F8").Value = 338
$ws.Range("F9").Value = 91
$ws.Range("F10").Value = 146
$ws.Range("F11").Value = 395
$ws.Range("F12").Value = 778
$ws.Range("F13").Value = 224
$ws.Range("F15").Value = 1116
$ws.Range("F16").Value = 298
$ws.Range("F17").Value = 651
$ws.Range("F18").Value = 1067
$ws.Range("F20").Value = 500
$ws.Range("F23").Value = 2857
$ws.Range("F28").Value = 7
$ws.Range("F29").Value = 219
$ws.Range("F31").Value = 2438
$ws.Range("F33").Value = 974
$ws.Range("F34").Value = 595
$ws.Range("F35").Value = 595
$ws.Range("F38").Value = 271
$ws.Range("F44").Value = 1056
$ws.Range("F46").Value = 37
$ws.Range("F49").Value = 29
